$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 and J1 (mirroring the style used by the other header cells)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data for columns I (I0) and J (IF), rows 2-24
$data = @(
    @(6, 7),
    @(5, 8),
    @(9, 9),
    @(11, 11),
    @(4, 6),
    @(6, 8),
    @(7, 7),
    @(7, 9),
    @(9, 9),
    @(7, 8),
    @(9, 9),
    @(7, 7),
    @(4, 5),
    @(7, 8),
    @(4, 6),
    @(8, 9),
    @(5, 6),
    @(1, 4),
    @(2, 5),
    @(6, 6),
    @(4, 4),
    @(3, 4),
    @(1, 1)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}
